$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

$holeIds = @(
    "LBU_05_18",
    "LBU_07_03",
    "LBU_05_08",
    "MHZ_08_05",
    "LBU_05_22",
    "LBU_05_21",
    "MHZ_12_02",
    "LBU_05_09",
    "LBU_05_06",
    "LBU_05_15",
    "LBU_05_24",
    "LBU_07_01",
    "LBU_05_23",
    "LBU_05_13",
    "MHZ_08_04",
    "LBU_05_29",
    "LBU_05_10",
    "LBU_01_3",
    "LBU_05_20",
    "LBU_05_27",
    "LBU_05_25",
    "LBU_05_02",
    "MHZ_12_01",
    "LBU_05_30",
    "LBU_05_17",
    "LBU_01_1",
    "LBU_02_3",
    "LBU_07_02",
    "LBU_05_26",
    "MHZ_08_02",
    "LBU_05_05",
    "LBU_05_04",
    "LBU_05_19",
    "MHZ_08_03",
    "LBU_01_2",
    "LBU_05_01",
    "MHZ_12_04"
)

$ws.Range("A1").Value = "hole_id"
$ws.Range("A2").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $holeIds[$i]
}
